$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update query text in B2, B3, B4 (append/adjust "order by ... limit" clauses) ---

$b2 = $ws.Range("B2").Text
$ws.Range("B2").Value = $b2 + "`n order By ss.study_subject_id ASC LIMIT 100"

$b3 = $ws.Range("B3").Text
$ws.Range("B3").Value = $b3 + "`n order By samp.sample_id ASC LIMIT 100"

$b4 = $ws.Range("B4").Text
$b4New = $b4.Replace("    order by f.file_name", "     order By f.file_name ASC LIMIT 100")
$ws.Range("B4").Value = $b4New

# --- Update selection (active cell moved from D4 to C4) ---
$ws.Range("C4").Select()

# --- Update row heights to reflect the extra wrapped line ---
$ws.Rows.Item(2).RowHeight = 331.2
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(4).RowHeight = 409.6
